# Update cryptos list rows 2-51 with the latest scraped values.
# Price (D) and Volume(1h) (E) cells must stay plain text exactly as
# scraped (e.g. "112.80", "29.376.83"), so we force a temporary text
# number format before writing, then clear formats again so no stray
# style survives the round-trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.376.83', '  -1.12%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.895.35', '  -1.48%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.002', '  +0.10%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '323.65', '  -3.60%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  +0.03%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4760', '  +1.43%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.4059', '  -1.94%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08025', '  -0.34%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.000', '  -1.67%  '),
    @(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '23.39', '  +4.37%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.910.26', '  -1.67%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.919', '  -1.81%  '),
    @(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.040', '  -2.30%  '),
    @(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '89.46', '  -0.71%  '),
    @(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  +0.11%  '),
    @(17, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06686', '  +1.41%  '),
    @(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001028', '  -0.66%  '),
    @(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.56', '  -1.82%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9999', '  +0.10%  '),
    @(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.380.86', '  -1.01%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.514', '  -0.95%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.68', '  +0.38%  '),
    @(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.149', '  -2.45%  '),
    @(25, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.069.92', '  -4.66%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '154.37', '  -1.52%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.76', '  -0.87%  '),
    @(28, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '6.011', '  +4.59%  '),
    @(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.083', '  -3.08%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '117.85', '  -0.05%  '),
    @(31, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.018', '  -3.95%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09460', '  -0.12%  '),
    @(33, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.384', '  -3.92%  '),
    @(34, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.529', '  +0.14%  '),
    @(35, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.350', '  -1.63%  '),
    @(36, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02246', '  -1.16%  '),
    @(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06033', '  -2.04%  '),
    @(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.169', '  -1.54%  '),
    @(39, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5853', '  -1.17%  '),
    @(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '7.804', '  -7.82%  '),
    @(41, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1836', '  -0.69%  '),
    @(42, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.09', '  -1.93%  '),
    @(43, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.287', '  +2.27%  '),
    @(44, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.411', '  +2.70%  '),
    @(45, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07712', '  +2.46%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.19', '  -0.67%  '),
    @(47, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5491', '  -1.95%  '),
    @(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.915', '  -1.57%  '),
    @(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '112.80', '  -0.10%  '),
    @(50, 'WOONetwork', 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo', '0.2966', '  -1.16%  '),
    @(51, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '43.54', '  -1.00%  ')

)

foreach ($row in $rows) {
    $r = $row[0]
    $coin = $row[1]
    $link = $row[2]
    $price = $row[3]
    $volume = $row[4]

    $ws.Range("B$r").Value = $coin
    $ws.Range("C$r").Value = $link

    $priceCell = $ws.Range("D$r")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.ClearFormats()

    $volCell = $ws.Range("E$r")
    $volCell.NumberFormat = "@"
    $volCell.Value = $volume
    $volCell.ClearFormats()
}

Write-Host "Updated $($rows.Count) rows"
